# HMPSTT (2018-08-01) 58_4 - correct district names to official names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose District (column G) value is "Bellary" / "Bellari" and should
# become the official name "Ballari (Bellary)".
$rows = @(4,5,6,7,8,10,12,15,19,22,23,24,25,27,30,31,32,37,38,39,40,43,44,47,53,54,55,56,60,61)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "Ballari (Bellary)"
}

# Rows 9 and 42 have a stray empty Address cell (column F) that should be
# removed entirely (the address text lives in column G for these rows).
$ws.Range("F9").ClearContents()
$ws.Range("F42").ClearContents()
